$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '246.93'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.463'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05617'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.469'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8037'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.050'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07255'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03176'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02962'
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09261'
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001663'
$ws.Range("E14").Value = '13BitForexTokenBF'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.226'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04686'
$ws.Range("E16").Value = '15CoinExTokenCET'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006277'
$ws.Range("E17").Value = '16TigerCashTCH'
$ws.Range("B18").Value = 'BitKan'
$ws.Range("C18").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.001050'
$ws.Range("E18").Value = '17BitKanKAN'
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.003813'
$ws.Range("E19").Value = '18HotbitTokenHTB'
$ws.Range("B20").Value = 'NitroEx'
$ws.Range("C20").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0001504'
$ws.Range("E20").Value = '19NitroExNTX'
$ws.Range("B21").Value = 'UpBots'
$ws.Range("C21").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0003608'
$ws.Range("E21").Value = '20UpBotsUBXT'
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.978'
$ws.Range("E22").Value = '21LEOLEO'
$ws.Range("B23").Value = 'GateToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.396'
$ws.Range("E23").Value = '22GateTokenGT'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.114'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("B25").Value = 'One'
$ws.Range("C25").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.01164'
$ws.Range("E25").Value = '24OneONEBestin24h'
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.3291'
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'
$ws.Range("B27").Value = 'ProBitToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1292'
$ws.Range("E27").Value = '26ProBitTokenPROB'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04156'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006897'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003508'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1040'
$ws.Range("E43").Value = '42BKEXTokenBKK'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01036'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005653'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.02620'
